# Transactions spreadsheet maintenance edit:
#  - Delete Transaction Button: remove the sample "bye" transaction row
#  - Fix None type error when adding data: remove the header row (the app
#    re-creates headers in code; leaving it in the data range caused the bug)
#  - Added Transaction Tab filtering function: append several new
#    transactions that were entered through the app after the fix

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Delete the "bye" row (row 2) first, then the header row (row 1), so the
# remaining transactions shift up to become rows 1-10.
$ws1.Rows.Item(2).Delete()
$ws1.Rows.Item(1).Delete()

# Append the new transactions entered after the fix.
$newRows = @(
    @("New Headphones", "US Bank Savings", "235.69", "01-01-2023", "Music!!!"),
    @("New Gun", "US Bank Savings", "350.68", "01-01-2023", "Bang, Bang"),
    @("Office Desk", "US Bank Checking", "608.79", "01-01-2023", ""),
    @("Pens", "Cash in wallet", "8.98", "01-01-2023", "Ran out of pens and needed more."),
    @("Powerwheel for Lucas", "US Bank Checking", "348.46", "01-01-2023", "Love that kid!!"),
    @("New computer", "US Bank Checking", "2300.98", "01-01-2023", "Very powerfull"),
    @("Tire for dirtbike", "US Bank Checking", "97.87", "01-01-2023", ""),
    @("New bib for tire", "US Bank Checking", "120.35", "01-01-2023", ""),
    @("Notebooks for work", "US Bank Checking", "18.59", "03-24-2023", "Gotta keep your notes!")
)

$startRow = 11
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    $ws1.Cells.Item($row, 1).Value = $data[0]
    $ws1.Cells.Item($row, 2).Value = $data[1]
    $ws1.Cells.Item($row, 3).Value = $data[2]
    $ws1.Cells.Item($row, 4).Value = $data[3]
    if ($data[4] -ne "") {
        $ws1.Cells.Item($row, 5).Value = $data[4]
    }
}

$ws1.Range("E19").Select()
